$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 4 (HAB Data) with the refined mapping
$ws.Range("A4").Value = "WQDataStations"
$ws.Range("B4").Value = "HABData"
$ws.Range("C4").Value = "Medium"
$ws.Range("D4").Value = "Multiple Sources"
$ws.Range("E4").Value = "This is Lake Erie data from stations that measure Chlorophyll and BlueGreen Algae levels."
$ws.Range("F4").Value = $null

# Row height for row 4 shrinks now that the link text is gone
$ws.Rows.Item(4).RowHeight = 45

# Update the sheet view/selection to match
$ws.Range("F3").Select()
